$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")

# --- Insert "large_wooden_farmhouse" as a new row before row 59 (naganuma) ---
# Clone formatting/formulas from row 62 ("wooden_house"), which is the template
# this new building variant is based on.
$ws.Rows.Item(62).Copy()
$ws.Rows.Item(59).Insert()

$ws.Range("A59").Value2 = "large_wooden_farmhouse"
$ws.Range("B59").Value2 = "large_wooden_farmhouse"
$ws.Range("C59").Value2 = 14
$ws.Range("F59").Value2 = "NAME_LARGE_WOODEN_FARMHOUSE"

# --- Insert "wooden_farmhouse" as a new row before the (now shifted) wooden_house row ---
# After the first insert, wooden_house moved from row 62 to row 63.
$ws.Rows.Item(63).Copy()
$ws.Rows.Item(63).Insert()

$ws.Range("A63").Value2 = "wooden_farmhouse"
$ws.Range("B63").Value2 = "wooden_farmhouse"
$ws.Range("C63").Value2 = 9
$ws.Range("F63").Value2 = "NAME_WOODEN_FARMHOUSE"

# --- Restore the sheet view (frozen pane / selection) to match the new layout ---
$ws.Activate()
$window = $excel.ActiveWindow
$window.SplitColumn = 1
$window.SplitRow = 1
$window.FreezePanes = $true

$ws.Range("B32").Select()
$ws.Range("A63").Select()
